# Actualización automática de grupos experimentales
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap / fill in the "Grupo_Experimental" assignments in column B
$ws.Range("B2").Value = "Sin SmartScore"
$ws.Range("B3").Value = "Con SmartScore"
$ws.Range("B7").Value = "Con SmartScore"
$ws.Range("B8").Value = "Con SmartScore"
$ws.Range("B9").Value = "Sin SmartScore"
$ws.Range("B10").Value = "Sin SmartScore"

# Row 10 SmartScore values were stored as text; convert them to real numbers
$ws.Range("I10").Value = 0.533
$ws.Range("L10").Value = 0.422
$ws.Range("O10").Value = 0.419
$ws.Range("R10").Value = 0.631
$ws.Range("U10").Value = 0.622
$ws.Range("X10").Value = 0.587
$ws.Range("AA10").Value = 0.762
$ws.Range("AD10").Value = 0.544
$ws.Range("AG10").Value = 0.497
